$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "datos actualizados" timestamp footer (last data row, column A)
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 08:52"

# 2) Countries table is sorted descending by "Casos totales" (column B).
#    Israel's totals were updated and it now jumps above Kuwait and
#    Republica Dominicana, so those three rows need to be re-sorted:
#    row 39 -> Israel (new data), row 40 -> Kuwait (old row-39 data),
#    row 41 -> Republica Dominicana (old row-40 data). Panama (row 42)
#    is untouched.
$ws.Range("A39").Value = "Israel"
$ws.Range("B39").Value = 64458
$ws.Range("C39").Value = 473
$ws.Range("D39").Value = 32109
$ws.Range("E39").Value = 31875
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 474

$ws.Range("A40").Value = "Kuwait"
$ws.Range("B40").Value = 64379
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 55057
$ws.Range("E40").Value = 8884
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 438

$ws.Range("A41").Value = "Republica Dominicana"
$ws.Range("B41").Value = 64156
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 30204
$ws.Range("E41").Value = 32869
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 1083

# 3) Straight data refreshes for the remaining countries (no re-sort needed,
#    each stays between its neighbours after the update).

# Row 53 - Armenia
$ws.Range("B53").Value = 37629
$ws.Range("C53").Value = 239
$ws.Range("D53").Value = 27357
$ws.Range("E53").Value = 9553
$ws.Range("G53").Value = 8
$ws.Range("H53").Value = 719

# Row 54 - Afganistan
$ws.Range("B54").Value = 36368
$ws.Range("C54").Value = 105
$ws.Range("D54").Value = 25358
$ws.Range("E54").Value = 9740
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 1270

# Row 75 - El Salvador
$ws.Range("D75").Value = 7793
$ws.Range("E75").Value = 6834

# Row 106 - Malaui
$ws.Range("B106").Value = 3709
$ws.Range("C106").Value = 45
$ws.Range("D106").Value = 1667
$ws.Range("E106").Value = 1939
$ws.Range("G106").Value = 4
$ws.Range("H106").Value = 103

# Row 139 - Letonia
$ws.Range("B139").Value = 1220
$ws.Range("C139").Value = 1
$ws.Range("D139").Value = 1052
$ws.Range("E139").Value = 137

# Row 143 - Georgia
$ws.Range("B143").Value = 1145
$ws.Range("C143").Value = 8
$ws.Range("D143").Value = 927
$ws.Range("E143").Value = 202

# Row 162 - Taiwan
$ws.Range("B162").Value = 467
$ws.Range("C162").Value = 5
$ws.Range("E162").Value = 20
